# HW Activity 4 edits
# ---------------------------------------------------------------
$d = $word.ActiveDocument

# 1. "Now, what you have all been waiting..." -> "Now, what you've all been waiting..."
$d.Content.Find.Execute(
    "you have all been waiting", $true, $false, $false, $false, $false,
    $true, 1, $false, "you" + [char]0x2019 + "ve all been waiting", 2) | Out-Null

# 2. "...manufacturer. These are different..." -> "...manufacturer. These drivers are different..."
$d.Content.Find.Execute(
    "the graphics card manufacturer. These are ", $true, $false, $false, $false, $false,
    $true, 1, $false, "the graphics card manufacturer. These drivers are ", 2) | Out-Null

# 3. "Install the graphics driver - again to save time, the driver is on the flash drive."
#    -> "Install the graphics driver - to save time, use the driver on the flash drive provided by your instructor."
$d.Content.Find.Execute(
    "Install the graphics driver - again to save time, the driver is on the flash drive.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Install the graphics driver - to save time, use the driver on the flash drive provided by your instructor.", 2) | Out-Null

# 4. "...Internet, e.g Google)..." -> "...Internet, e.g. Google)..."
$d.Content.Find.Execute(
    "the Internet, e.g Google", $true, $false, $false, $false, $false,
    $true, 1, $false, "the Internet, e.g. Google", 2) | Out-Null

# 4b. Insert bold "motherboard " before the bold "manufacturer's web[site]" run.
$anchor = $d.Content
$anchor.Find.Execute(
    "is to get any required drivers from the ", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0) | Out-Null
$insPos = $anchor.End
$ins = $d.Range($insPos, $insPos)
$ins.InsertBefore("motherboard ")
$boldNew = $d.Range($insPos, $insPos + 12)
$boldNew.Font.Bold = 1

# 4c. "tw). " -> "tw). Write the manufacturer and model number for your motherboard below."
$d.Content.Find.Execute(
    "tw). ", $true, $false, $false, $false, $false,
    $true, 1, $false, "tw). Write the manufacturer and model number for your motherboard below.", 2) | Out-Null

# 5. Drop the stale lastRenderedPageBreak by touching the run it lives on.
$d.Content.Find.Execute(
    "Using the ", $true, $false, $false, $false, $false,
    $true, 1, $false, "Using the ", 2) | Out-Null

# 6. BIOS sentence rewrite.
$d.Content.Find.Execute(
    "s. You can also often get an updated BIOS should one be available to provide added support for newer hardware.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "s. Oftentimes, you can also find an updated BIOS that adds support for newer hardware.", 2) | Out-Null

# 7. "missing ones" -> "missing drivers" + move the _GoBack bookmark there.
$d.Content.Find.Execute(
    "be additional categories or missing ones if the current operating system includes them",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "be additional categories or missing drivers if the current operating system includes them", 2) | Out-Null

$gb = $d.Content
$gb.Find.Execute("missing drivers", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$gbPoint = $d.Range($gb.End, $gb.End)
$d.Bookmarks.Add("_GoBack", $gbPoint)
